$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '65.958.58'
$ws.Range('E2').Value2 = '  -1.05%  '
$ws.Range('D3').Value2 = '3.519.77'
$ws.Range('E3').Value2 = '  +0.39%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value2 = '1.00'
$c.Style = 'Normal'
$ws.Range('E4').Value2 = '  -0.14%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value2 = '577.51'
$c.Style = 'Normal'
$ws.Range('E5').Value2 = '  +4.51%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value2 = '178.75'
$c.Style = 'Normal'
$ws.Range('E6').Value2 = '  -6.07%  '
$ws.Range('E7').Value2 = '  +4.85%  '
$ws.Range('E8').Value2 = '  -0.06%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value2 = '0.635'
$c.Style = 'Normal'
$ws.Range('E9').Value2 = '  +0.10%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value2 = '0.157'
$c.Style = 'Normal'
$ws.Range('E10').Value2 = '  +5.09%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value2 = '55.33'
$c.Style = 'Normal'
$ws.Range('E11').Value2 = '  -0.34%  '
$ws.Range('E12').Value2 = '  +1.99%  '
$ws.Range('E13').Value2 = '  -1.89%  '
$ws.Range('D14').Value2 = '4.079.67'
$ws.Range('E14').Value2 = '  +0.13%  '
$ws.Range('D15').Value2 = '3.518.80'
$ws.Range('E15').Value2 = '  +0.44%  '
$ws.Range('E16').Value2 = '  +0.22%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value2 = '18.45'
$c.Style = 'Normal'
$ws.Range('E17').Value2 = '  +0.98%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value2 = '12.13'
$c.Style = 'Normal'
$ws.Range('E18').Value2 = '  +2.69%  '
$ws.Range('D19').Value2 = '65.917.34'
$ws.Range('E19').Value2 = '  -1.26%  '
$ws.Range('E20').Value2 = '  +1.41%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value2 = '415.72'
$c.Style = 'Normal'
$ws.Range('E21').Value2 = '  +1.11%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value2 = '4.24'
$c.Style = 'Normal'
$ws.Range('E22').Value2 = '  +7.77%  '
$ws.Range('B23').Value2 = 'Toncoin'
$ws.Range('C23').Value2 = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value2 = '4.34'
$c.Style = 'Normal'
$ws.Range('E23').Value2 = '  +3.51%  '
$ws.Range('B24').Value2 = 'Litecoin'
$ws.Range('C24').Value2 = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value2 = '85.89'
$c.Style = 'Normal'
$ws.Range('E24').Value2 = '  +0.67%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value2 = '13.07'
$c.Style = 'Normal'
$ws.Range('E25').Value2 = '  +9.83%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value2 = '11.00'
$c.Style = 'Normal'
$ws.Range('E26').Value2 = '  -1.22%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value2 = '2.87'
$c.Style = 'Normal'
$ws.Range('E27').Value2 = '  -2.27%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value2 = '9.07'
$c.Style = 'Normal'
$ws.Range('E28').Value2 = '  +2.58%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value2 = '30.40'
$c.Style = 'Normal'
$ws.Range('E29').Value2 = '  +0.29%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value2 = '625.04'
$c.Style = 'Normal'
$ws.Range('E30').Value2 = '  -4.75%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value2 = '6.48'
$c.Style = 'Normal'
$ws.Range('E31').Value2 = '  -3.22%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value2 = '11.70'
$c.Style = 'Normal'
$ws.Range('E32').Value2 = '  -0.51%  '
$ws.Range('E33').Value2 = '  -0.36%  '
$ws.Range('B34').Value2 = 'Kaspa'
$ws.Range('C34').Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value2 = '0.155'
$c.Style = 'Normal'
$ws.Range('E34').Value2 = '  +13.32%  '
$ws.Range('B35').Value2 = 'OKB'
$ws.Range('C35').Value2 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value2 = '59.70'
$c.Style = 'Normal'
$ws.Range('E35').Value2 = '  -0.30%  '
$ws.Range('D36').Value2 = '0.0₃0804'
$ws.Range('E36').Value2 = '  -0.27%  '
$ws.Range('E37').Value2 = '  +0.29%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value2 = '37.36'
$c.Style = 'Normal'
$ws.Range('E38').Value2 = '  -3.77%  '
$ws.Range('B39').Value2 = 'Maker'
$ws.Range('C39').Value2 = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value2 = '3.302.18'
$ws.Range('E39').Value2 = '  +10.15%  '
$ws.Range('B40').Value2 = 'TheGraph'
$ws.Range('C40').Value2 = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value2 = '0.381'
$c.Style = 'Normal'
$ws.Range('E40').Value2 = '  -2.85%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value2 = '3.46'
$c.Style = 'Normal'
$ws.Range('E41').Value2 = '  +3.64%  '
$ws.Range('E42').Value2 = '  -0.21%  '
$ws.Range('E43').Value2 = '  -3.67%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value2 = '0.0418'
$c.Style = 'Normal'
$ws.Range('E44').Value2 = '  +0.32%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value2 = '2.51'
$c.Style = 'Normal'
$ws.Range('E45').Value2 = '  -4.76%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value2 = '3.24'
$c.Style = 'Normal'
$ws.Range('E46').Value2 = '  -3.54%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value2 = '2.73'
$c.Style = 'Normal'
$ws.Range('E47').Value2 = '  +0.27%  '
$ws.Range('E48').Value2 = '  +2.18%  '
$ws.Range('B49').Value2 = 'Monero'
$ws.Range('C49').Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value2 = '140.55'
$c.Style = 'Normal'
$ws.Range('E49').Value2 = '  +1.15%  '
$ws.Range('B50').Value2 = 'THORChain'
$ws.Range('C50').Value2 = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value2 = '8.60'
$c.Style = 'Normal'
$ws.Range('E50').Value2 = '  -3.98%  '
$ws.Range('E51').Value2 = '  -4.80%  '
